$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 data update: new trade search result (Rustington job) replaces
# the old Angmering/Littlehampton/Dan Hunt row.
$ws.Range("A5").Value = 230761
$ws.Range("D5").Value = "David Clarke"
$ws.Range("B5").Value = "Rustington"
$ws.Range("C5").Value = "Rustington, West Sussex"

# The pasted location text in C5 carries its own (web-sourced) font
# (Arial, #333333) rather than the sheet's default Calibri/theme color.
$tempStyle = $wb.Styles.Add("__PastedWebText")
$tempStyle.Font.Color = 3355443
$tempStyle.Font.Name = "Arial"
$ws.Range("C5").Style = "__PastedWebText"
$wb.Styles.Item("__PastedWebText").Delete()

# Selection left on C5 after the edit.
$ws.Range("C5").Select()
